# ccdi dataset icdc changes
# Adds a new "CoreAdditional" worksheet (dbGaP core/additional dataset metadata)
# after the existing "DatasetPage" sheet, and updates view/selection state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. New worksheet, positioned right after DatasetPage
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "CoreAdditional"

# ---------------------------------------------------------------------------
# 2. Column widths (bestFit widths carried over from the authoring session)
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 33.085
$ws2.Columns.Item(2).ColumnWidth = 51.926656
$ws2.Columns.Item(3).ColumnWidth = 10.423344
$ws2.Columns.Item(4).ColumnWidth = 21.259969
$ws2.Columns.Item(5).ColumnWidth = 27.926656
$ws2.Columns.Item(6).ColumnWidth = 15.926656
$ws2.Columns.Item(7).ColumnWidth = 15.926656
$ws2.Columns.Item(8).ColumnWidth = 10.923344
$ws2.Columns.Item(9).ColumnWidth = 8.926656
$ws2.Columns.Item(10).ColumnWidth = 35.256656

# ---------------------------------------------------------------------------
# 3. Header row
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "suffixUrl"
$ws2.Range("B1").Value = "DatasetName"
$ws2.Range("C1").Value = "numOfCases"
$ws2.Range("D1").Value = "caseSex"
$ws2.Range("E1").Value = "caseAge"
$ws2.Range("F1").Value = "caseRace"
$ws2.Range("G1").Value = "caseEthn"
$ws2.Range("H1").Value = "noOfSamples"
$ws2.Range("I1").Value = "dbgapID"
$ws2.Range("J1").Value = "grant"

# ---------------------------------------------------------------------------
# 4. Row 2 - dbGaP-phs001928 dataset
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "/dataset/dbGaP-phs001928"
$ws2.Range("B2").Value = "A Comprehensive Genomic Study of Pediatric Malignancy"
$ws2.Range("C2").Value = 267
$ws2.Range("D2").Value = "Not Reported (267)"
$ws2.Range("E2").Value = "Pediatric and Young Adult (<40 years) (267)"
$ws2.Range("F2").Value = "Not Reported (267)"
$ws2.Range("G2").Value = "Not Reported (267)"
$ws2.Range("H2").Value = 401
$ws2.Range("I2").Value = "phs001928"
$ws2.Range("J2").Value = "Intramural Research Program of the Center for Cancer Research, NCI"

# ---------------------------------------------------------------------------
# 5. Row 3 - TARGET-ALL Phase 2 dataset
# ---------------------------------------------------------------------------
$ws2.Range("A3").Value = "/dataset/TARGET-ALL%20Phase%202"
$ws2.Range("B3").Value = "Acute Lymphoblastic Leukemia (ALL) Expansion Phase 2"
$ws2.Range("C3").Value = 791
$ws2.Range("D3").Value = "Female (361); Male (430)"
$ws2.Range("E3").Value = "0 to 4 years (248); 5 to 9 years (126); 10 to 14 years (252); 15 to 19 years (152); 20 to 24 years (11); 25 to 29 years (2)"
$ws2.Range("F3").Value = "American Indian or Alaska Native (3); Asian (39); Black or African American (58); Native Hawaiian or Other Pacific Islander (2); White (579); Unknown (110)"
$ws2.Range("G3").Value = "Hispanic or Latino (212); Not Hispanic or Latino (546); Unknown (33)"
$ws2.Range("I3").Value = "phs000464"
$ws2.Range("J3").Value = "261200800001E-12-0-40`nTherapeutically Applicable Research to Generate Effective Treatments (TARGET) `nHHSN261200800001E`nNCI Contract `nU10CA180886`nCOG NCTN Network Group Operations Center"

# ---------------------------------------------------------------------------
# 6. Number formats - whole used range is Text ("@" -> numFmtId 49)
# ---------------------------------------------------------------------------
$ws2.Range("A1:J3").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 7. Wrap text for the long-form columns/cells
# ---------------------------------------------------------------------------
$ws2.Range("J1:J3").WrapText = $true
$ws2.Range("E3:G3").WrapText = $true

# ---------------------------------------------------------------------------
# 8. A2 looks like the hyperlinked suffixUrl cell on DatasetPage
# ---------------------------------------------------------------------------
$ws2.Range("A2").Style = "Hyperlink"
$ws2.Range("A2").NumberFormat = "@"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://datacatalog-qa.ccdi.cancer.gov/dataset/dbGaP-phs001928", [System.Type]::Missing, [System.Type]::Missing, "https://datacatalog-qa.ccdi.cancer.gov/dataset/dbGaP-phs001928") | Out-Null

# ---------------------------------------------------------------------------
# 9. E2 uses a small Lato font (matches the source site's rendering)
# ---------------------------------------------------------------------------
$ws2.Range("E2").Font.Name = "Lato"
$ws2.Range("E2").Font.Size = 8
$ws2.Range("E2").Font.Color = 2696481

# ---------------------------------------------------------------------------
# 10. Row heights for the wrapped rows
# ---------------------------------------------------------------------------
$ws2.Rows.Item(2).RowHeight = 29.5
$ws2.Rows.Item(3).RowHeight = 130.5

# ---------------------------------------------------------------------------
# 11. Page setup + hyperlink-navigation selection on the new sheet
# ---------------------------------------------------------------------------
$ws2.PageSetup.Orientation = 1
$ws2.Range("B14").Select()

# ---------------------------------------------------------------------------
# 12. DatasetPage view state moves off the active tab
# ---------------------------------------------------------------------------
$ws1.Range("B8").Select()

# ---------------------------------------------------------------------------
# 13. Activate the new sheet last so it becomes the active/visible tab
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B14").Select()
